$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "36.396.39"
$ws.Range("E2").Value = "  -0.60%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.930.77"
$ws.Range("E3").Value = "  -4.05%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "242.22"
$ws.Range("E5").Value = "  -2.04%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -4.02%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "56.65"
$ws.Range("E8").Value = "  -10.15%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -5.48%  "

# Row 10 - OKB
Set-TextValue $ws.Range("D10") "55.27"
$ws.Range("E10").Value = "  -3.36%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0828"
$ws.Range("E11").Value = "  +4.48%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -0.59%  "

# Row 13 - Polygon
Set-TextValue $ws.Range("D13") "0.815"
$ws.Range("E13").Value = "  -7.93%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "21.34"
$ws.Range("E14").Value = "  -7.01%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.214.00"
$ws.Range("E15").Value = "  -4.12%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "13.34"
$ws.Range("E16").Value = "  -6.10%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  -6.56%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "1.930.48"
$ws.Range("E18").Value = "  -4.20%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "36.311.82"
$ws.Range("E19").Value = "  -0.62%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "69.16"
$ws.Range("E20").Value = "  -3.90%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  -2.06%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "228.02"
$ws.Range("E22").Value = "  -4.43%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -7.44%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.08%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -3.09%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -2.46%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.25"
$ws.Range("E27").Value = "  -8.15%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "162.52"
$ws.Range("E28").Value = "  +1.95%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "19.23"
$ws.Range("E29").Value = "  -4.77%  "

# Row 30 - Kaspa
$ws.Range("E30").Value = "  -7.12%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -3.62%  "

# Row 32 - ImmutableX
Set-TextValue $ws.Range("D32") "1.14"
$ws.Range("E32").Value = "  -3.27%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -7.74%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0626"
$ws.Range("E34").Value = "  -1.27%  "

# Row 35 - InternetComputer(DFINITY)
$ws.Range("E35").Value = "  -5.43%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  +0.06%  "

# Row 37 - THORChain
$ws.Range("E37").Value = "  -7.31%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  -2.89%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -9.25%  "

# Row 40 - RenderToken
$ws.Range("E40").Value = "  -9.10%  "

# Row 41 - Cronos
Set-TextValue $ws.Range("D41") "0.0968"
$ws.Range("E41").Value = "  -3.88%  "

# Row 42 - HuobiToken
$ws.Range("E42").Value = "  -1.96%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -7.54%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -5.14%  "

# Row 45 - InjectiveProtocol
Set-TextValue $ws.Range("D45") "15.67"
$ws.Range("E45").Value = "  -6.50%  "

# Row 46 & 47 swap: Maker <-> ARBITRUM
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D46") "1.03"
$ws.Range("E46").Value = "  -8.80%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D47") "1.339.78"
$ws.Range("E47").Value = "  -1.77%  "

# Row 48 - Aave
Set-TextValue $ws.Range("D48") "87.00"
$ws.Range("E48").Value = "  -9.52%  "

# Row 49 - FraxShare
Set-TextValue $ws.Range("D49") "7.19"
$ws.Range("E49").Value = "  -6.72%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  -2.88%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "45.36"
$ws.Range("E51").Value = "  +1.44%  "
